$wb = $excel.ActiveWorkbook

# Worksheets involved
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
# Both data rows share the same timestamp value, update both.
$wsOverview.Range("G2").Value = "2016-09-06 05:37:39"
$wsOverview.Range("G3").Value = "2016-09-06 05:37:39"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-09-06 05:37:29"
$wsZhCn.Range("H3").Value = "2016-09-06 05:37:29"

# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-09-06 05:38:11"
$wsZhCn.Range("K3").Value = "2016-09-06 05:38:11"

# --- de-de sheet ---
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Correspond Handoff Datetime column (H) - shares value with Overview's G column
$wsDeDe.Range("H2").Value = "2016-09-06 05:37:39"
$wsDeDe.Range("H3").Value = "2016-09-06 05:37:39"

# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-09-06 05:38:27"
$wsDeDe.Range("K3").Value = "2016-09-06 05:38:27"
